# Updated cryptos list on Tue Dec 12 05:09:13 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string value into a cell while forcing text storage
# (so numeric-looking strings like "251.98" don't get silently converted
# to a Number by Excel's type inference) and without leaving any lasting
# style/format change behind on the cell.
function Set-TextValue {
    param($range, [string]$text)
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $savedStyle
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "41.545.17"
Set-TextValue $ws.Range("E2") "  -1.70%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.222.55"
Set-TextValue $ws.Range("E3") "  -1.06%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  -0.05%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "251.98"
Set-TextValue $ws.Range("E5") "  +8.30%  "

# Row 6 - XRP
Set-TextValue $ws.Range("E6") "  -1.40%  "

# Row 7 - Solana
Set-TextValue $ws.Range("D7") "70.63"
Set-TextValue $ws.Range("E7") "  +0.21%  "

# Row 8 - USDC
Set-TextValue $ws.Range("E8") "  +0.00%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.563"
Set-TextValue $ws.Range("E9") "  -0.51%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "42.11"
Set-TextValue $ws.Range("E10") "  +15.81%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.0960"
Set-TextValue $ws.Range("E11") "  -4.08%  "

# Row 12 - OKB
Set-TextValue $ws.Range("D12") "58.74"
Set-TextValue $ws.Range("E12") "  +0.25%  "

# Row 13 - TRON
Set-TextValue $ws.Range("E13") "  +0.24%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "6.94"
Set-TextValue $ws.Range("E14") "  +0.61%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "2.544.83"
Set-TextValue $ws.Range("E15") "  -1.21%  "

# Row 16 - Chainlink
Set-TextValue $ws.Range("D16") "14.88"
Set-TextValue $ws.Range("E16") "  -1.36%  "

# Row 17 - Polygon
Set-TextValue $ws.Range("D17") "0.853"
Set-TextValue $ws.Range("E17") "  -1.77%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "2.218.19"
Set-TextValue $ws.Range("E18") "  -1.21%  "

# Row 19 - WrappedBTC
Set-TextValue $ws.Range("D19") "41.449.60"
Set-TextValue $ws.Range("E19") "  -1.56%  "

# Row 20 - ShibaInu
Set-TextValue $ws.Range("D20") "0.0₃0963"
Set-TextValue $ws.Range("E20") "  -2.56%  "

# Row 21 - was Uniswap, now Litecoin
Set-TextValue $ws.Range("B21") "Litecoin"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D21") "72.85"
Set-TextValue $ws.Range("E21") "  -1.01%  "

# Row 22 - was Litecoin, now Uniswap
Set-TextValue $ws.Range("B22") "Uniswap"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D22") "6.18"
Set-TextValue $ws.Range("E22") "  -1.42%  "

# Row 23 - ImmutableX
Set-TextValue $ws.Range("E23") "  +11.46%  "

# Row 24 - BitcoinCash
Set-TextValue $ws.Range("D24") "233.96"
Set-TextValue $ws.Range("E24") "  -1.73%  "

# Row 25 - Dai
Set-TextValue $ws.Range("E25") "  +0.03%  "

# Row 26 - WEMIXToken
Set-TextValue $ws.Range("D26") "3.77"
Set-TextValue $ws.Range("E26") "  +3.32%  "

# Row 27 - PancakeSwap
Set-TextValue $ws.Range("D27") "2.51"
Set-TextValue $ws.Range("E27") "  +6.28%  "

# Row 28 - Cosmos
Set-TextValue $ws.Range("D28") "10.22"
Set-TextValue $ws.Range("E28") "  +0.68%  "

# Row 29 - Toncoin
Set-TextValue $ws.Range("E29") "  +1.44%  "

# Row 30 - Monero
Set-TextValue $ws.Range("D30") "171.11"
Set-TextValue $ws.Range("E30") "  +1.56%  "

# Row 31 - EthereumClassic
Set-TextValue $ws.Range("E31") "  -0.95%  "

# Row 32 - Kaspa
Set-TextValue $ws.Range("E32") "  +0.09%  "

# Row 33 - Stellar
Set-TextValue $ws.Range("D33") "0.124"
Set-TextValue $ws.Range("E33") "  -2.42%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("E34") "  +2.25%  "

# Row 35 - Hedera
Set-TextValue $ws.Range("D35") "0.0716"
Set-TextValue $ws.Range("E35") "  -0.71%  "

# Row 36 - InjectiveProtocol
Set-TextValue $ws.Range("D36") "26.28"
Set-TextValue $ws.Range("E36") "  +17.76%  "

# Row 37 - Filecoin
Set-TextValue $ws.Range("D37") "4.66"
Set-TextValue $ws.Range("E37") "  -3.31%  "

# Row 38 - RenderToken
Set-TextValue $ws.Range("D38") "3.96"
Set-TextValue $ws.Range("E38") "  +9.25%  "

# Row 39 - VeChain
Set-TextValue $ws.Range("E39") "  +5.87%  "

# Row 40 - LidoDAOToken
Set-TextValue $ws.Range("E40") "  +1.18%  "

# Row 41 - MultiversX
Set-TextValue $ws.Range("D41") "69.17"
Set-TextValue $ws.Range("E41") "  +2.56%  "

# Row 42 - THORChain
Set-TextValue $ws.Range("D42") "5.99"
Set-TextValue $ws.Range("E42") "  -1.99%  "

# Row 43 - Celestia
Set-TextValue $ws.Range("D43") "12.01"
Set-TextValue $ws.Range("E43") "  +16.80%  "

# Row 44 - was Algorand, now FTXToken
Set-TextValue $ws.Range("B44") "FTXToken"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D44") "5.07"
Set-TextValue $ws.Range("E44") "  +1.13%  "

# Row 45 - was FTXToken, now Algorand
Set-TextValue $ws.Range("B45") "Algorand"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D45") "0.208"
Set-TextValue $ws.Range("E45") "  +9.97%  "

# Row 46 - SynthetixNetwork
Set-TextValue $ws.Range("D46") "4.77"
Set-TextValue $ws.Range("E46") "  +7.55%  "

# Row 47 - FraxShare
Set-TextValue $ws.Range("D47") "8.74"
Set-TextValue $ws.Range("E47") "  -3.66%  "

# Row 48 - Cronos
Set-TextValue $ws.Range("E48") "  +0.29%  "

# Row 49 - BinanceUSD
Set-TextValue $ws.Range("E49") "  -0.11%  "

# Row 50 - ARBITRUM
Set-TextValue $ws.Range("E50") "  +5.54%  "

# Row 51 - TrustWalletToken
Set-TextValue $ws.Range("E51") "  +0.66%  "
